$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the mAP1/mAP2/mAP3 column headers to MT / MT+T / MT+T+SS
$ws.Range("B1").Value = "MT"
$ws.Range("C1").Value = "MT+T"
$ws.Range("D1").Value = "MT+T+SS"

# Move the active selection from H40 to D1
$ws.Range("D1").Select()

# Configure page setup: A4 paper, portrait orientation
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
